$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 3 - L.Fournette
$rushing.Cells.Item(3, 3).Value = 103   # C3
$rushing.Cells.Item(3, 4).Value = 32    # D3
$rushing.Cells.Item(3, 5).Value = 17    # E3
$rushing.Cells.Item(3, 6).Value = 39    # F3

# Row 4 - R.Jones
$rushing.Cells.Item(4, 3).Value = 39    # C4

# Row 5 - G.Bernard
$rushing.Cells.Item(5, 3).Value = 3     # C5

# Row 8 - J.Darden
$rushing.Cells.Item(8, 4).Value = 3     # D8

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - R.Jones
$receiving.Cells.Item(2, 3).Value = 67   # C2
$receiving.Cells.Item(2, 4).Value = 56   # D2
$receiving.Cells.Item(2, 7).Value = 15   # G2
$receiving.Cells.Item(2, 8).Value = 11   # H2

# Row 3 - G.Bernard
$receiving.Cells.Item(3, 3).Value = 6    # C3
$receiving.Cells.Item(3, 4).Value = 4    # D3

# Row 4 - K.Vaughn
$receiving.Cells.Item(4, 3).Value = 26   # C4
$receiving.Cells.Item(4, 4).Value = 21   # D4

# Row 5 - C.Godwin
$receiving.Cells.Item(5, 3).Value = 65   # C5
$receiving.Cells.Item(5, 4).Value = 43   # D5
$receiving.Cells.Item(5, 5).Value = 28   # E5
$receiving.Cells.Item(5, 6).Value = 15   # F5

# Row 6 - J.Darden
$receiving.Cells.Item(6, 3).Value = 85   # C6
$receiving.Cells.Item(6, 4).Value = 69   # D6
$receiving.Cells.Item(6, 5).Value = 21   # E6
$receiving.Cells.Item(6, 6).Value = 13   # F6
$receiving.Cells.Item(6, 7).Value = 23   # G6
$receiving.Cells.Item(6, 8).Value = 18   # H6

# Row 8 - T.Johnson
$receiving.Cells.Item(8, 3).Value = 25   # C8
$receiving.Cells.Item(8, 4).Value = 16   # D8
$receiving.Cells.Item(8, 7).Value = 5    # G8
$receiving.Cells.Item(8, 8).Value = 2    # H8

# Row 11 - B.Perriman
$receiving.Cells.Item(11, 3).Value = 6   # C11
$receiving.Cells.Item(11, 4).Value = 3   # D11

# Row 12 - R.Gronkowski
$receiving.Cells.Item(12, 3).Value = 38  # C12
$receiving.Cells.Item(12, 4).Value = 29  # D12
$receiving.Cells.Item(12, 5).Value = 16  # E12
$receiving.Cells.Item(12, 6).Value = 13  # F12
$receiving.Cells.Item(12, 7).Value = 10  # G12
$receiving.Cells.Item(12, 8).Value = 7   # H12

# Row 14 - C.Brate
$receiving.Cells.Item(14, 3).Value = 34  # C14
$receiving.Cells.Item(14, 4).Value = 20  # D14
$receiving.Cells.Item(14, 7).Value = 16  # G14
$receiving.Cells.Item(14, 8).Value = 9   # H14
